# Des Scheduled Flights vs actual.xlsx
# Append 26 new days of data (2021-08-09 .. 2021-09-03) below the existing table,
# which previously ended at row 490 (2021-08-08).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 490
$firstNewRow = 491
$lastNewRow = 516

# 1) Stamp out the new rows by copying the formatting (styles/number formats)
#    of the last existing row (A:D) down into each new row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Range("A" + $lastRow + ":D" + $lastRow).Copy($ws.Range("A" + $r + ":D" + $r))
}

# 2) Write the real DateTime / Scheduled flights / Tracked flights values.
$ws.Range("A491").Value = "2021-08-09"
$ws.Range("B491").Value = 73
$ws.Range("C491").Value = 65
$ws.Range("A492").Value = "2021-08-10"
$ws.Range("B492").Value = 69
$ws.Range("C492").Value = 61
$ws.Range("A493").Value = "2021-08-11"
$ws.Range("B493").Value = 77
$ws.Range("C493").Value = 70
$ws.Range("A494").Value = "2021-08-12"
$ws.Range("B494").Value = 81
$ws.Range("C494").Value = 80
$ws.Range("A495").Value = "2021-08-13"
$ws.Range("B495").Value = 81
$ws.Range("C495").Value = 74
$ws.Range("A496").Value = "2021-08-14"
$ws.Range("B496").Value = 64
$ws.Range("C496").Value = 62
$ws.Range("A497").Value = "2021-08-15"
$ws.Range("B497").Value = 75
$ws.Range("C497").Value = 71
$ws.Range("A498").Value = "2021-08-16"
$ws.Range("B498").Value = 72
$ws.Range("C498").Value = 70
$ws.Range("A499").Value = "2021-08-17"
$ws.Range("B499").Value = 68
$ws.Range("C499").Value = 66
$ws.Range("A500").Value = "2021-08-18"
$ws.Range("B500").Value = 73
$ws.Range("C500").Value = 66
$ws.Range("A501").Value = "2021-08-19"
$ws.Range("B501").Value = 85
$ws.Range("C501").Value = 83
$ws.Range("A502").Value = "2021-08-20"
$ws.Range("B502").Value = 67
$ws.Range("C502").Value = 63
$ws.Range("A503").Value = "2021-08-21"
$ws.Range("B503").Value = 64
$ws.Range("C503").Value = 61
$ws.Range("A504").Value = "2021-08-22"
$ws.Range("B504").Value = 67
$ws.Range("C504").Value = 63
$ws.Range("A505").Value = "2021-08-23"
$ws.Range("B505").Value = 67
$ws.Range("C505").Value = 67
$ws.Range("A506").Value = "2021-08-24"
$ws.Range("B506").Value = 74
$ws.Range("C506").Value = 65
$ws.Range("A507").Value = "2021-08-25"
$ws.Range("B507").Value = 67
$ws.Range("C507").Value = 64
$ws.Range("A508").Value = "2021-08-26"
$ws.Range("B508").Value = 73
$ws.Range("C508").Value = 70
$ws.Range("A509").Value = "2021-08-27"
$ws.Range("B509").Value = 63
$ws.Range("C509").Value = 61
$ws.Range("A510").Value = "2021-08-28"
$ws.Range("B510").Value = 61
$ws.Range("C510").Value = 58
$ws.Range("A511").Value = "2021-08-29"
$ws.Range("B511").Value = 65
$ws.Range("C511").Value = 64
$ws.Range("A512").Value = "2021-08-30"
$ws.Range("B512").Value = 77
$ws.Range("C512").Value = 69
$ws.Range("A513").Value = "2021-08-31"
$ws.Range("B513").Value = 69
$ws.Range("C513").Value = 67
$ws.Range("A514").Value = "2021-09-01"
$ws.Range("B514").Value = 73
$ws.Range("C514").Value = 66
$ws.Range("A515").Value = "2021-09-02"
$ws.Range("B515").Value = 80
$ws.Range("C515").Value = 74
$ws.Range("A516").Value = "2021-09-03"
$ws.Range("B516").Value = 71
$ws.Range("C516").Value = 70

# 3) Fill the "% on time" formula (Tracked/Scheduled) down column D for the new rows.
#    Done as two fill actions, matching the two ranges in the saved workbook.
$ws.Range("D491:D514").Formula = "=C491/B491"
$ws.Range("D515:D516").Formula = "=C515/B515"

# 4) Restore the selection that was active when the workbook was last saved.
$ws.Range("F514").Select()
